$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 37-40: copy formatting from row 36 (column A style) before filling values ---
$ws.Range("A36").Copy()
$ws.Range("A37:A40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column A sequence numbers for the 4 new rows ---
$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(39, 1).Value = 37
$ws.Cells.Item(40, 1).Value = 38

# --- Column B: new "Buying Opportunity" ticker list, rows 2-40 ---
$ws.Cells.Item(2, 2).Value = "NSE:ASIANENE"
$ws.Cells.Item(3, 2).Value = "NSE:AUROPHARMA"
$ws.Cells.Item(4, 2).Value = "NSE:BANKBARODA"
$ws.Cells.Item(5, 2).Value = "NSE:BANKINDIA"
$ws.Cells.Item(6, 2).Value = "NSE:BAYERCROP"
$ws.Cells.Item(7, 2).Value = "NSE:BHARATRAS"
$ws.Cells.Item(8, 2).Value = "NSE:BHEL"
$ws.Cells.Item(9, 2).Value = "NSE:BIRLACORPN"
$ws.Cells.Item(10, 2).Value = "NSE:CALSOFT"
$ws.Cells.Item(11, 2).Value = "NSE:CASTROLIND"
$ws.Cells.Item(12, 2).Value = "NSE:DEVIT"
$ws.Cells.Item(13, 2).Value = "NSE:DIVISLAB"
$ws.Cells.Item(14, 2).Value = "NSE:DIVOPPBEES"
$ws.Cells.Item(15, 2).Value = "NSE:GLAND"
$ws.Cells.Item(16, 2).Value = "NSE:HDFCNEXT50"
$ws.Cells.Item(17, 2).Value = "NSE:HGS"
$ws.Cells.Item(18, 2).Value = "NSE:HPAL"
$ws.Cells.Item(19, 2).Value = "NSE:INDIGO"
$ws.Cells.Item(20, 2).Value = "NSE:INDOCO"
$ws.Cells.Item(21, 2).Value = "NSE:INGERRAND"
$ws.Cells.Item(22, 2).Value = "NSE:JIOFIN"
$ws.Cells.Item(23, 2).Value = "NSE:JMA"
$ws.Cells.Item(24, 2).Value = "NSE:KINGFA"
$ws.Cells.Item(25, 2).Value = "NSE:LICMFGOLD"
$ws.Cells.Item(26, 2).Value = "NSE:LUMAXIND"
$ws.Cells.Item(27, 2).Value = "NSE:MEGASTAR"
$ws.Cells.Item(28, 2).Value = "NSE:MOHEALTH"
$ws.Cells.Item(29, 2).Value = "NSE:MOM100"
$ws.Cells.Item(30, 2).Value = "NSE:MOQUALITY"
$ws.Cells.Item(31, 2).Value = "NSE:NEXT50"
$ws.Cells.Item(32, 2).Value = "NSE:NLCINDIA"
$ws.Cells.Item(33, 2).Value = "NSE:PATINTLOG"
$ws.Cells.Item(34, 2).Value = "NSE:PHARMABEES"
$ws.Cells.Item(35, 2).Value = "NSE:PITTIENG"
$ws.Cells.Item(36, 2).Value = "NSE:POWERMECH"
$ws.Cells.Item(37, 2).Value = "NSE:PRAJIND"
$ws.Cells.Item(38, 2).Value = "NSE:PROZONER"
$ws.Cells.Item(39, 2).Value = "NSE:PTL"
$ws.Cells.Item(40, 2).Value = "NSE:RESPONIND"

# --- Column C: "support Zone" tickers shift up for rows 2-5; rows 6-40 cleared ---
$ws.Cells.Item(2, 3).Value = "NSE:APOLLO"
$ws.Cells.Item(3, 3).Value = "NSE:CYIENT"
$ws.Cells.Item(4, 3).Value = "NSE:GRSE"
$ws.Cells.Item(5, 3).Value = "NSE:RRKABEL"
$ws.Cells.Item(6, 3).Value = ""
$ws.Cells.Item(7, 3).Value = ""
$ws.Cells.Item(8, 3).Value = ""
$ws.Cells.Item(9, 3).Value = ""
$ws.Cells.Item(10, 3).Value = ""
$ws.Cells.Item(11, 3).Value = ""
$ws.Cells.Item(12, 3).Value = ""
$ws.Cells.Item(13, 3).Value = ""
$ws.Cells.Item(14, 3).Value = ""
$ws.Cells.Item(15, 3).Value = ""
$ws.Cells.Item(16, 3).Value = ""
$ws.Cells.Item(17, 3).Value = ""
$ws.Cells.Item(18, 3).Value = ""
$ws.Cells.Item(19, 3).Value = ""
$ws.Cells.Item(20, 3).Value = ""
$ws.Cells.Item(21, 3).Value = ""
$ws.Cells.Item(22, 3).Value = ""
$ws.Cells.Item(23, 3).Value = ""
$ws.Cells.Item(24, 3).Value = ""
$ws.Cells.Item(25, 3).Value = ""
$ws.Cells.Item(26, 3).Value = ""
$ws.Cells.Item(27, 3).Value = ""
$ws.Cells.Item(28, 3).Value = ""
$ws.Cells.Item(29, 3).Value = ""
$ws.Cells.Item(30, 3).Value = ""
$ws.Cells.Item(31, 3).Value = ""
$ws.Cells.Item(32, 3).Value = ""
$ws.Cells.Item(33, 3).Value = ""
$ws.Cells.Item(34, 3).Value = ""
$ws.Cells.Item(35, 3).Value = ""
$ws.Cells.Item(36, 3).Value = ""
$ws.Cells.Item(37, 3).Value = ""
$ws.Cells.Item(38, 3).Value = ""
$ws.Cells.Item(39, 3).Value = ""
$ws.Cells.Item(40, 3).Value = ""

# --- Column E: "Short buildup" tickers removed for rows 2-5 ---
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(3, 5).Value = ""
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(5, 5).Value = ""
